$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 32,2
$arr[0,0] = -0.24779303665832941
$arr[0,1] = 0.24759424728270574
$arr[1,0] = -0.17858012829818115
$arr[1,1] = 0.17807940251682464
$arr[2,0] = -0.12837497082954386
$arr[2,1] = 0.12804905231882557
$arr[3,0] = -0.12004905233703766
$arr[3,1] = 0.11949294530657362
$arr[4,0] = -0.1164929453172725
$arr[4,1] = 0.11459113033694823
$arr[5,0] = -0.010335881732736141
$arr[5,1] = 0.010210989266926873
$arr[6,0] = -0.0067446028468141073
$arr[6,1] = 0.00672034925122178
$arr[7,0] = 0.0032796507226602323
$arr[7,1] = -0.0032979195366960035
$arr[8,0] = 0.0052979195246911637
$arr[8,1] = -0.0053076089017185524
$arr[9,0] = 0.0073076088898673675
$arr[9,1] = -0.0073075947812508701
$arr[10,0] = 0.010307594767648531
$arr[10,1] = -0.010308310748856186
$arr[11,0] = 0.013808310734420903
$arr[11,1] = -0.013823714788694375
$arr[12,0] = 0.0173237147748182
$arr[12,1] = -0.017341826852090847
$arr[13,0] = 0.025341826830586278
$arr[13,1] = -0.025379680575587038
$arr[14,0] = 0.026379680566772201
$arr[14,1] = -0.026431268506288852
$arr[15,0] = -0.0060329649972867827
$arr[15,1] = 0.0060035311199606589
$arr[16,0] = -0.0040035311302011323
$arr[16,1] = 0.0039999999861910496
$arr[17,0] = -0.063014118995603496
$arr[17,1] = 0.062920712466421946
$arr[18,0] = -0.012092108468504836
$arr[18,1] = 0.012017240249716288
$arr[19,0] = -0.008017240258018532
$arr[19,1] = 0.0080057315608090818
$arr[20,0] = -0.0040057315691983675
$arr[20,1] = 0.0039999999915458773
$arr[21,0] = -0.045704431697247472
$arr[21,1] = 0.04549312564030572
$arr[22,0] = -0.040493125652532491
$arr[22,1] = 0.040097948642451087
$arr[23,0] = -0.020097948682002986
$arr[23,1] = 0.019999999959969372
$arr[24,0] = -0.097220800435295018
$arr[24,1] = 0.097096811824496143
$arr[25,0] = -0.094596811837796224
$arr[25,1] = 0.094437815344743825
$arr[26,0] = -0.091937815358732689
$arr[26,1] = 0.091000745414713258
$arr[27,0] = -0.089000745430853456
$arr[27,1] = 0.088359054574217666
$arr[28,0] = -0.081359054601407799
$arr[28,1] = 0.081172666506382996
$arr[29,0] = -0.021172666627708114
$arr[29,1] = 0.021024653762442114
$arr[30,0] = -0.014024653791578245
$arr[30,1] = 0.014001692059588677
$arr[31,0] = -0.0040016920941194911
$arr[31,1] = 0.0039999999759281479

$ws.Range("A1:B32").Value = $arr
